$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct small text issues in the sandwich descriptions (column B)
$ws.Range("B12").Value = "Pepperoni, Salami, 2x Provolone, Jalapeño Peppers, Lettuce, Tomatoes, Red Onions, MVP Parmesan Vinaigrette"
$ws.Range("B13").Value = "Ham, Capicola, BelGioioso Fresh Mozzarella, Spinach, Tomatoes, Red Onions, Banana peppers"
$ws.Range("B14").Value = "Pepperoni, Meatballs, BelGioioso Fresh Mozzarella, Parmesan, Toasted"
$ws.Range("B19").Value = "Rotisserie-Style Chicken, 2x PepperJack, Lettuce, Tomatoes, Red Onions, Baja Chipotle"
$ws.Range("B23").Value = "Turkey, Bacon, BelGioioso Fresh Mozzarella, Toasted, Smashed Avocado, Spinach, Tomatoes, Red Onions, Mayo"

# Move the active cell selection from D7 to B5
[void]$ws.Range("B5").Select()
